$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 112438847
$ws.Range("B2").Value = 95707
$ws.Range("Q2").Value = 502260
$ws.Range("R2").Value = 6543183

# Row 3
$ws.Range("A3").Value = 112438846
$ws.Range("B3").Value = 95704
$ws.Range("E3").Value = 221946
$ws.Range("F3").Value = "Mattlummer"
$ws.Range("G3").Value = "Lycopodium clavatum"
$ws.Range("Q3").Value = 502269
$ws.Range("R3").Value = 6543231

# Row 4
$ws.Range("A4").Value = 112438848
$ws.Range("B4").Value = 95707
$ws.Range("E4").Value = 221941
$ws.Range("F4").Value = "Plattlummer"
$ws.Range("G4").Value = "Lycopodium complanatum"
$ws.Range("Q4").Value = 502199
$ws.Range("R4").Value = 6543178

# Row 5
$ws.Range("B5").Value = 95704

# Row 6
$ws.Range("A6").Value = 112438849
$ws.Range("B6").Value = 95704
$ws.Range("Q6").Value = 502192
$ws.Range("R6").Value = 6543228
